$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear existing content, keep A1's format around long enough to clone it ---
$ws.Cells.ClearContents()

# A1 still carries the original bold/centered/bordered header style (s="1").
# Clone that format onto the new header cells (B1:E1) and onto column A
# for the data rows (A2:A9), then drop the stale A1 cell entirely.
$ws.Range("A1").Copy()
$ws.Range("B1:E1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A2:A9").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A1").Clear()

# --- Header row (row 1) ---
$ws.Range("B1").Value = "lang_code"
$ws.Range("C1").Value = "code"
$ws.Range("D1").Value = "descr"
$ws.Range("E1").Value = "is_active"

# --- Data rows (2-9) ---
$data = @(
    @(0, "eng", "txt",  "Text File",     $true),
    @(1, "eng", "xml",  "XML File",      $true),
    @(2, "eng", "json", "Json File",     $true),
    @(3, "fra", "txt",  "Fichier texte", $true),
    @(4, "fra", "xml",  "Fichier XML",   $true),
    @(5, "fra", "json", "Fichier Json",  $true),
    @(6, "eng", "html", "html file",     $true),
    @(7, "fra", "html", "Fichier html",  $true)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}
